$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds text-typed values that look numeric (e.g. "1.00",
# "36.665.66" - the latter using "." as a thousands separator rather than a decimal
# point). Assigning such strings directly to .Value lets Excel auto-detect them as
# numbers and silently reformat/round them, so for those cells we briefly force Text
# number format, assign the literal text, then restore the default (Normal) style.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.665.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.964.23'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.84%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.620'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.14'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.374'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0817'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.40%  '
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.252.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.830'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.963.79'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.557.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '229.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("E25").Value = '  +3.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.143'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +13.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.24'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +1.77%  '
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.72'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0618'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("E35").Value = '  +6.47%  '
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.93%  '
$ws.Range("E38").Value = '  -3.31%  '
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0989'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("E41").Value = '  +1.43%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.17'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0212'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.360.96'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.70%  '
$ws.Range("E46").Value = '  +0.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.79'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.15'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("E49").Value = '  +0.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.142.82'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.90'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.54%  '
